# Add three new backlog rows (47-49) to the "Developmnet PB" sheet, mirroring
# the styling of the existing rows (column A: wrap-text label, column B:
# status-colour fill, column C: status text), and update the sheet's
# scroll position / active selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Developmnet PB")
$ws.Activate()

# --- New backlog items ------------------------------------------------
$ws.Range("A47").Value = "Improve view pages layout for equipment"
$ws.Range("A48").Value = "Add numeric handling for equipment attributes (important for no. of patient handlers)"
$ws.Range("A49").Value = "Check anywhere that there is a completeattribute object used in the model that the field variables are valid"

# Column A uses the wrap-text label style used throughout the sheet
$ws.Range("A47:A49").WrapText = $true

# Column B status-colour fill cells (no values, just formatting)
$ws.Range("B47").Interior.Color = 49407   # orange - matches row 46 (B46 s=6)
$ws.Range("B48").Interior.Color = 49407   # orange - matches row 46 (B46 s=6)
$ws.Range("B49").Interior.Color = 255     # red - matches row 5 (B5 s=5)

# Column C status text - "Not Started" for rows 47 and 48 (row 49 has none)
$ws.Range("C47").Value = "Not Started"
$ws.Range("C48").Value = "Not Started"

# --- Update the view ----------------------------------------------------
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("C39").Select()
